$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at position 102 - this pushes the former rows
# 102-115 down to 103-116 (and carries the date-format style of column D
# down with them automatically).
$ws.Rows.Item(102).Insert()

# Populate the newly inserted row 102 with the new weekly record.
$ws.Range("A102").Value = 4
$ws.Range("B102").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C102").Value = "Los Lagos"
$ws.Range("D102").Value = 44504
$ws.Range("E102").Value = 10
$ws.Range("F102").Value = "Fruta"
$ws.Range("G102").Value = 100101
$ws.Range("H102").Value = "Berries"
$ws.Range("I102").Value = 100112025
$ws.Range("J102").Value = "Frutilla"
$ws.Range("K102").Value = "Sin especificar"
$ws.Range("L102").Value = "Primera"
$ws.Range("M102").Value = 200
$ws.Range("N102").Value = 9000
$ws.Range("O102").Value = 10000
$ws.Range("P102").Value = 9500
$ws.Range("Q102").Value = "`$/bandeja 7 kilos"
$ws.Range("R102").Value = "Provincia de Melipilla"
$ws.Range("S102").Value = 1357
$ws.Range("T102").Value = 7
